# Re-order the per-profile statistics rows (2-5) so they read, top to bottom,
# PHYSICS, MEDICINE, MATHEMATICS, LINGUISTICS (previously LINGUISTICS,
# MATHEMATICS, PHYSICS, MEDICINE), carrying each profile's figures and
# university list along with its label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "PHYSICS";     B = 4.539999961853027; C = 8.0; D = 2.0; E = "Московский Выдуманный Университет;Московский Придуманный Институт;" },
    @{ A = "MEDICINE";    B = 4.329999923706055; C = 3.0; D = 3.0; E = "Московский Государственный Медицинский Университет;Тамбовский Университет Медицины;Самарский Медицинский Институт;" },
    @{ A = "MATHEMATICS"; B = 0.0;               C = 0.0; D = 1.0; E = "Казанский Университет Вычислений;" },
    @{ A = "LINGUISTICS"; B = 0.0;               C = 0.0; D = 1.0; E = "Воронежский Литературно-Переводческий Университет;" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}
